$wb = $excel.ActiveWorkbook

# --- Worksheets ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1. Update status text "Ready for handoff" -> "In Translation" ---
# Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# Per-language detail sheets: Status column C, rows 2-3
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- 2. Narrow the status columns ---
# Target stored column width (OOXML) is 13.4101845877511 characters.
# Excel's ColumnWidth COM property differs from the stored XML width by the
# standard 5-pixel padding (5/6 of a character at the default font), so we
# back that out before assigning.
$newColumnWidth = 13.4101845877511 - (5 / 6)

# Overview sheet: columns E and F
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# zh-cn / de-de detail sheets: column C
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
